$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G14").Value = 3.1
$ws.Range("I14").Value = 2.1
$ws.Range("Q14").Value = 2.7
$ws.Range("T14").Value = 9.75
$ws.Range("U14").Value = 16.5
$ws.Range("V14").Value = 11
$ws.Range("W14").Value = 40
$ws.Range("X14").Value = 27
$ws.Range("Y14").Value = 35
$ws.Range("AA14").Value = 6.5
$ws.Range("AD14").Value = 7.9
$ws.Range("AE14").Value = 10.5
$ws.Range("AF14").Value = 8.75
$ws.Range("AG14").Value = 20
$ws.Range("AH14").Value = 16.5
$ws.Range("AI14").Value = 26
